# Remove the false-start data rows (old rows 2 and 3), shifting the
# remaining rows (old 4, old 5) up into rows 2 and 3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("2:3").Delete()

# Reselect the now-current data rows, matching the author's last selection.
$ws.Range("A2:XFD3").Select()
